$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 40 (hunk 0)
$ws.Range("H40").Value = 1669383.4
$ws.Range("I40").Value = 3334666.8
$ws.Range("K40").Value = 3334666.8
$ws.Range("M40").Value = -3334491.8
# Row 42 (hunk 1)
$ws.Range("H42").Value = 391.6
$ws.Range("I42").Value = 108
$ws.Range("J42").Value = 462.5
$ws.Range("K42").Value = 324
$ws.Range("L42").Value = 1387.5
$ws.Range("M42").Value = -94
$ws.Range("N42").Value = -1847.5
# Row 58 (hunk 2)
$ws.Range("H58").Value = 465.5
$ws.Range("I58").Value = 322.375
$ws.Range("J58").Value = 751.75
$ws.Range("K58").Value = 967.125
$ws.Range("L58").Value = 2255.25
$ws.Range("M58").Value = -817.125
$ws.Range("N58").Value = -2555.25
# Row 132 (hunk 3)
$ws.Range("H132").Value = 2235.5757
$ws.Range("I132").Value = 2146.0667
$ws.Range("K132").Value = 6438.2001
$ws.Range("M132").Value = -3908.2001

$ws = $wb.Worksheets.Item("ARM")
# Row 61 (hunk 4)
$ws.Range("H61").Value = 11908669
$ws.Range("I61").Value = 16131814
$ws.Range("K61").Value = 16131814
$ws.Range("M61").Value = -16131602
# Row 74 (hunk 5)
$ws.Range("H74").Value = 33372258
$ws.Range("I74").Value = 45506704
$ws.Range("J74").Value = 2538.375
$ws.Range("K74").Value = 45506704
$ws.Range("L74").Value = 2538.375
$ws.Range("M74").Value = -45505830
$ws.Range("N74").Value = -4286.375
# Row 77 (hunk 6)
$ws.Range("H77").Value = 33372258
$ws.Range("I77").Value = 45506704
$ws.Range("J77").Value = 2538.375
$ws.Range("K77").Value = 227533520
$ws.Range("L77").Value = 12691.875
$ws.Range("M77").Value = -227529152
$ws.Range("N77").Value = -21427.875
# Row 110 (hunk 7)
$ws.Range("H110").Value = 31974.584
$ws.Range("I110").Value = 45713.125
$ws.Range("J110").Value = 4497.5
$ws.Range("K110").Value = 45713.125
$ws.Range("L110").Value = 4497.5
$ws.Range("M110").Value = -43668.125
$ws.Range("N110").Value = -8587.5
# Row 132 (hunk 8)
$ws.Range("H132").Value = 43550604
$ws.Range("I132").Value = 14172.111
$ws.Range("J132").Value = 200281760
$ws.Range("K132").Value = 42516.333
$ws.Range("L132").Value = 600845280
$ws.Range("M132").Value = -39986.333
$ws.Range("N132").Value = -600850340
# Row 136 (hunk 9)
$ws.Range("H136").Value = 11908669
$ws.Range("I136").Value = 16131814
$ws.Range("K136").Value = 48395442
$ws.Range("M136").Value = -48392892

$ws = $wb.Worksheets.Item("BSM")
# Row 61 (hunk 10)
$ws.Range("H61").Value = 105000
$ws.Range("J61").Value = 105000
$ws.Range("L61").Value = 105000
$ws.Range("N61").Value = -105626
# Row 134 (hunk 11)
$ws.Range("H134").Value = 2440.1667
$ws.Range("I134").Value = 2239.1177
$ws.Range("K134").Value = 6717.353099999999
$ws.Range("M134").Value = -4182.353099999999

$ws = $wb.Worksheets.Item("CRP")
# Row 132 (hunk 12)
$ws.Range("H132").Value = 45900.957
$ws.Range("I132").Value = 62430.21
$ws.Range("J132").Value = 3942.077
$ws.Range("K132").Value = 187290.63
$ws.Range("L132").Value = 11826.231
$ws.Range("M132").Value = -184760.63
$ws.Range("N132").Value = -16886.231
# Row 134 (hunk 13)
$ws.Range("H134").Value = 850
$ws.Range("I134").Value = 850
$ws.Range("K134").Value = 2550
$ws.Range("M134").Value = -15

$ws = $wb.Worksheets.Item("CUL")
# Row 64 (hunk 14)
$ws.Range("H64").Value = 9031
$ws.Range("I64").Value = 2932
$ws.Range("K64").Value = 8796
$ws.Range("M64").Value = -8526
# Row 67 (hunk 15)
$ws.Range("H67").Value = 9031
$ws.Range("I67").Value = 2932
$ws.Range("K67").Value = 8796
$ws.Range("M67").Value = -7860
# Row 87 (hunk 16)
$ws.Range("H87").Value = 1000
$ws.Range("I87").Value = 1000
$ws.Range("K87").Value = 3000
$ws.Range("M87").Value = -1752
# Row 90 (hunk 17)
$ws.Range("H90").Value = 1000
$ws.Range("I90").Value = 1000
$ws.Range("K90").Value = 9000
$ws.Range("M90").Value = -2760
# Row 114 (hunk 18)
$ws.Range("H114").Value = 2511718.5
$ws.Range("I114").Value = 3334041.2
$ws.Range("J114").Value = 44750
$ws.Range("K114").Value = 10002123.6
$ws.Range("L114").Value = 134250
$ws.Range("M114").Value = -9998869.600000001
$ws.Range("N114").Value = -140758
# Row 122 (hunk 19)
$ws.Range("H122").Value = 1268.6428
$ws.Range("I122").Value = 412.6
$ws.Range("J122").Value = 1744.2222
$ws.Range("K122").Value = 3713.4
$ws.Range("L122").Value = 15697.9998
$ws.Range("M122").Value = -1263.4
$ws.Range("N122").Value = -20597.9998
# Row 131 (hunk 20)
$ws.Range("H131").Value = 1618.4706
$ws.Range("I131").Value = 1288
$ws.Range("J131").Value = 1675.4482
$ws.Range("K131").Value = 3864
$ws.Range("L131").Value = 5026.3446
$ws.Range("M131").Value = 1176
$ws.Range("N131").Value = -15106.3446
# Row 132 (hunk 21)
$ws.Range("H132").Value = 2714
$ws.Range("I132").Value = 2195
$ws.Range("J132").Value = 3129.2
$ws.Range("K132").Value = 19755
$ws.Range("L132").Value = 28162.8
$ws.Range("M132").Value = -17225
$ws.Range("N132").Value = -33222.8
# Row 137 (hunk 22)
$ws.Range("H137").Value = 5111
$ws.Range("J137").Value = 7200
$ws.Range("L137").Value = 21600
$ws.Range("N137").Value = -31800

$ws = $wb.Worksheets.Item("GSM")
# Row 2 (hunk 23)
$ws.Range("H2").Value = 2272970
$ws.Range("I2").Value = 3571530
$ws.Range("J2").Value = 489.875
$ws.Range("K2").Value = 3571530
$ws.Range("L2").Value = 489.875
$ws.Range("M2").Value = -3571417
$ws.Range("N2").Value = -715.875
# Row 131 (hunk 24)
$ws.Range("H131").Value = 88999.664
$ws.Range("J131").Value = 88999.664
$ws.Range("L131").Value = 88999.664
$ws.Range("N131").Value = -99079.664
# Row 132 (hunk 25)
$ws.Range("H132").Value = 2257.7144
$ws.Range("I132").Value = 2077.889
$ws.Range("J132").Value = 3336.6667
$ws.Range("K132").Value = 6233.667
$ws.Range("L132").Value = 10010.0001
$ws.Range("M132").Value = -3703.667
$ws.Range("N132").Value = -15070.0001

$ws = $wb.Worksheets.Item("LTW")
# Row 22 (hunk 26)
$ws.Range("H22").Value = 2710
$ws.Range("I22").Value = 625
$ws.Range("J22").Value = 4100
$ws.Range("K22").Value = 625
$ws.Range("L22").Value = 4100
$ws.Range("M22").Value = -330
$ws.Range("N22").Value = -4690
# Row 27 (hunk 27)
$ws.Range("H27").Value = 2710
$ws.Range("I27").Value = 625
$ws.Range("J27").Value = 4100
$ws.Range("K27").Value = 625
$ws.Range("L27").Value = 4100
$ws.Range("M27").Value = -518
$ws.Range("N27").Value = -4314
# Row 46 (hunk 28)
$ws.Range("H46").Value = 1248.5
$ws.Range("I46").Value = 667.34784
$ws.Range("J46").Value = 3921.8
$ws.Range("K46").Value = 667.34784
$ws.Range("L46").Value = 3921.8
$ws.Range("M46").Value = -479.34784
$ws.Range("N46").Value = -4297.8
# Row 68 (hunk 29)
$ws.Range("H68").Value = 3000
$ws.Range("I68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("M68").ClearContents()  # cell removed entirely in target
# Row 71 (hunk 30)
$ws.Range("H71").Value = 3000
$ws.Range("I71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("M71").ClearContents()  # cell removed entirely in target
# Row 132 (hunk 31)
$ws.Range("H132").Value = 12118.615
$ws.Range("I132").Value = 10452.639
$ws.Range("K132").Value = 31357.917
$ws.Range("M132").Value = -28827.917
# Row 136 (hunk 32)
$ws.Range("H136").Value = 2714.9375
$ws.Range("I136").Value = 699.875
$ws.Range("K136").Value = 2099.625
$ws.Range("M136").Value = 450.375

$ws = $wb.Worksheets.Item("WVR")
# Row 18 (hunk 33)
$ws.Range("H18").Value = 6750
$ws.Range("J18").Value = 6750
$ws.Range("L18").Value = 6750
$ws.Range("N18").Value = -7096
# Row 74 (hunk 34)
$ws.Range("H74").Value = 43992.25
$ws.Range("J74").Value = 43992.25
$ws.Range("L74").Value = 43992.25
$ws.Range("N74").Value = -45864.25
# Row 77 (hunk 35)
$ws.Range("H77").Value = 43992.25
$ws.Range("J77").Value = 43992.25
$ws.Range("L77").Value = 131976.75
$ws.Range("N77").Value = -141336.75
# Row 100 (hunk 36)
$ws.Range("H100").Value = 1594
$ws.Range("I100").Value = 2240
$ws.Range("K100").Value = 4480
$ws.Range("M100").Value = -3939
# Row 103 (hunk 37)
$ws.Range("H103").Value = 16846
$ws.Range("J103").Value = 16846
$ws.Range("L103").Value = 16846
$ws.Range("N103").Value = -19190
# Row 107 (hunk 38)
$ws.Range("H107").Value = 329.7
$ws.Range("I107").Value = 342.42856
$ws.Range("J107").Value = 300
$ws.Range("K107").Value = 1027.28568
$ws.Range("L107").Value = 900
$ws.Range("M107").Value = 892.71432
$ws.Range("N107").Value = -4740
# Row 113 (hunk 39)
$ws.Range("H113").Value = 626.35
$ws.Range("J113").Value = 856.2
$ws.Range("L113").Value = 2568.6
$ws.Range("N113").Value = -6908.6
# Row 122 (hunk 40)
$ws.Range("H122").Value = 23226.701
$ws.Range("I122").Value = 26589.324
$ws.Range("J122").Value = 4011.7144
$ws.Range("K122").Value = 79767.97200000001
$ws.Range("L122").Value = 12035.1432
$ws.Range("M122").Value = -77317.97200000001
$ws.Range("N122").Value = -16935.1432
